# Daily attendance processing - 2025-10-11 22:41:25
# Swap the order of the first two "Recorded By" entries (column G) for every
# data row on the "Session Analysis Results" sheet. Rows whose G value has
# only a single recorder (no comma) or is empty are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val -notmatch ",") {
        continue
    }

    $parts = $val -split ",\s*"
    if ($parts.Count -lt 2) {
        continue
    }

    $first = $parts[0]
    $second = $parts[1]

    if ($parts.Count -gt 2) {
        $rest = $parts[2..($parts.Count - 1)]
    } else {
        $rest = @()
    }

    $newParts = @($second, $first) + $rest
    $newVal = [string]::Join(", ", $newParts)

    $cell.Value2 = $newVal
}
